# HW: Add AC104B, AC105A, AC106A to BOM
#
# Appends three new BOM sub-assembly sections (AC104B, AC105A, AC106A) to
# the bottom of the "Sheet1" BOM table, following the same layout used by
# the existing sections (e.g. AC104A at rows 107-110): a bold/underlined
# section-header row holding the assembly name in column A, followed by
# one data row per component, then a blank (skipped) row before the next
# section.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats - clones the section-header look (bold font + bottom
# border on A:I) from an existing header row without disturbing the
# shared style table.
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# AC104B (rows 112-116)
# ---------------------------------------------------------------------
$ws.Range("A107:I107").Copy()
$ws.Range("A112:I112").PasteSpecial($xlPasteFormats)
$ws.Range("A112").Value = "AC104B"

$ws.Range("A113").Value = "CABLE 6COND 28AWG BLK 153M (0.15M)"
$ws.Range("B113").Value = "N/A"
$ws.Range("C113").Value = "Tensility"
$ws.Range("D113").Value = "30-00510"
$ws.Range("E113").Value = 1
$ws.Range("F113").Value = "Digi-Key"
$ws.Range("G113").Value = "T1347-153-ND"
$ws.Range("H113").Value = "N/A"
$ws.Range("I113").Value = "N/A"

$ws.Range("A114").Value = "CONN JACK FMALE 6POS SOLDER CUP"
$ws.Range("B114").Value = "N/A"
$ws.Range("C114").Value = "Hirose Electric Co Ltd"
$ws.Range("D114").Value = "HR10A-7J-6S(73)"
$ws.Range("E114").Value = 1
$ws.Range("F114").Value = "Digi-Key"
$ws.Range("G114").Value = "HR1601-ND"
$ws.Range("H114").Value = "Mouser"
$ws.Range("I114").Value = "798-HR10A-7J-6S73"

$ws.Range("A115").Value = "CONN MOD PLUG 8P8C UNSHIELDED"
$ws.Range("B115").Value = "N/A"
$ws.Range("C115").Value = "Stewart Connector"
$ws.Range("D115").Value = "940-SP-3088"
$ws.Range("E115").Value = 1
$ws.Range("F115").Value = "Digi-Key"
$ws.Range("G115").Value = "380-1192-ND"
$ws.Range("H115").Value = "Mouser"
$ws.Range("I115").Value = "530-940-SP-3088"

$ws.Range("A116").Value = "DIODE SCHOTTKY 40V SGL DO35"
$ws.Range("B116").Value = "N/A"
$ws.Range("C116").Value = "Vishay Semiconductor Diodes"
$ws.Range("D116").Value = "SD103A-TAP"
$ws.Range("E116").Value = 1
$ws.Range("F116").Value = "Digi-Key"
$ws.Range("G116").Value = "SD103A-TAPGICT-ND"
$ws.Range("H116").Value = "Mouser"
$ws.Range("I116").Value = "78-SD103A-TAP"

# row 117 intentionally left blank (separator, matches existing sections)

# ---------------------------------------------------------------------
# AC105A (rows 118-120)
# ---------------------------------------------------------------------
$ws.Range("A107:I107").Copy()
$ws.Range("A118:I118").PasteSpecial($xlPasteFormats)
$ws.Range("A118").Value = "AC105A"

$ws.Range("A119").Value = "BOX ABS BLACK 1.38""L X 1.38""W"
$ws.Range("B119").Value = "N/A"
$ws.Range("C119").Value = "Hammond Manufacturing"
$ws.Range("D119").Value = "1551MBK"
$ws.Range("E119").Value = 1
$ws.Range("F119").Value = "Digi-Key"
$ws.Range("G119").Value = "HM993-ND"
$ws.Range("H119").Value = "Mouser"
$ws.Range("I119").Value = "546-1551MBK"

$ws.Range("A120").Value = "CONN RCPT FMALE 6POS SOLDER CUP"
$ws.Range("B120").Value = "N/A"
$ws.Range("C120").Value = "Hirose Electric Co Ltd"
$ws.Range("D120").Value = "HR10A-7R-6S(73)"
$ws.Range("E120").Value = 4
$ws.Range("F120").Value = "Digi-Key"
$ws.Range("G120").Value = "HR1594-ND"
$ws.Range("H120").Value = "Mouser"
$ws.Range("I120").Value = "798-HR10A-7R-6S73"

# row 121 intentionally left blank (separator, matches existing sections)

# ---------------------------------------------------------------------
# AC106A (rows 122-126)
# ---------------------------------------------------------------------
$ws.Range("A107:I107").Copy()
$ws.Range("A122:I122").PasteSpecial($xlPasteFormats)
$ws.Range("A122").Value = "AC106A"

$ws.Range("A123").Value = "CABLE 6COND 28AWG BLK 153M (0.15M)"
$ws.Range("B123").Value = "N/A"
$ws.Range("C123").Value = "Tensility"
$ws.Range("D123").Value = "30-00510"
$ws.Range("E123").Value = 1
$ws.Range("F123").Value = "Digi-Key"
$ws.Range("G123").Value = "T1347-153-ND"
$ws.Range("H123").Value = "N/A"
$ws.Range("I123").Value = "N/A"

$ws.Range("A124").Value = "CONN JACK FMALE 6POS SOLDER CUP"
$ws.Range("B124").Value = "N/A"
$ws.Range("C124").Value = "Hirose Electric Co Ltd"
$ws.Range("D124").Value = "HR10A-7J-6S(73)"
$ws.Range("E124").Value = 1
$ws.Range("F124").Value = "Digi-Key"
$ws.Range("G124").Value = "HR1601-ND"
$ws.Range("H124").Value = "Mouser"
$ws.Range("I124").Value = "798-HR10A-7J-6S73"

$ws.Range("A125").Value = "Kenwood Portable Accessory Connector"
$ws.Range("B125").Value = "N/A"
$ws.Range("C125").Value = "N/A"
$ws.Range("D125").Value = "N/A"
$ws.Range("E125").Value = 1
$ws.Range("F125").Value = "N/A"
$ws.Range("G125").Value = "N/A"
$ws.Range("H125").Value = "N/A"
$ws.Range("I125").Value = "N/A"

$ws.Range("A126").Value = "DIODE SCHOTTKY 40V SGL DO35"
$ws.Range("B126").Value = "N/A"
$ws.Range("C126").Value = "Vishay Semiconductor Diodes"
$ws.Range("D126").Value = "SD103A-TAP"
$ws.Range("E126").Value = 1
$ws.Range("F126").Value = "Digi-Key"
$ws.Range("G126").Value = "SD103A-TAPGICT-ND"
$ws.Range("H126").Value = "Mouser"
$ws.Range("I126").Value = "78-SD103A-TAP"

Write-Output "Added AC104B, AC105A, AC106A sections (rows 112-126)"
